$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (pushes existing rows 7-16 down to 8-17),
# mirroring the author's "add a new weekly record" edit.
$ws.Rows(7).Insert()

# Populate the new row 7 with the new record. Most columns replicate the
# constant values shared by every data row in this sheet; only the date
# (D) and volume (J) differ for the newly-added entry.
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 'Macroferia Regional de Talca'
$ws.Range("C7").Value = 'Maule'
$ws.Range("D7").Value = 44701
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 100112040
$ws.Range("G7").Value = 'Cilantro'
$ws.Range("H7").Value = 'Sin especificar'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 7000
$ws.Range("N7").Value = '$/caja 36 atados'
$ws.Range("O7").Value = 'Región del Maule'
$ws.Range("P7").Value = 194
$ws.Range("Q7").Value = 36
$ws.Range("R7").Value = 'Hortaliza'
